$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("database")

# "database" -> "threats_database" (commit: "thay doi bang database -> trustedDatabase" / "them bang threatsDatabase")
$ws.Name = "threats_database"

# New threat record on the renamed sheet
$ws.Range("A2").Value = "192.168.10.1"
$ws.Range("B2").Value = "1.1.1.1.1.1"

# New sheet "trusted_database" placed right after "threats_database",
# carrying the header + the original row that used to live in "database"
$ws2 = $wb.Worksheets.Add($null, $ws, $null, $null)
$ws2.Name = "trusted_database"

$ws2.Range("A1").Value = "src ip"
$ws2.Range("B1").Value = "src mac"
$ws2.Range("A2").Value = "1.2.3.4"
$ws2.Range("B2").Value = "10.10.10.10.10.10"

# Keep the first sheet ("threats_database") as the active tab, same as before the edit
$ws.Activate()

Write-Host "Sheets:"
foreach ($s in $wb.Worksheets) {
    Write-Host " -" $s.Name()
}
Write-Host "Active sheet:" $wb.ActiveSheet.Name()
